# Change footer fontsize to 11pt
#
# The "Footer" paragraph style (and its linked "Footer Char" character
# style) currently render at 10.5pt (w:sz 21 half-points). Bump both to
# 11pt (w:sz 22 half-points) so footer text throughout the document picks
# up the new size.

$d = $word.ActiveDocument

$footerStyle = $d.Styles("Footer")
$footerStyle.Font.Size = 11

$footerCharStyle = $d.Styles("Footer Char")
$footerCharStyle.Font.Size = 11
